$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1 & 2) Split "The document must be a .docx Microsoft Word file" into three
#    runs so "docx" is wrapped in proofErr spellStart/spellEnd markers, i.e.:
#       "The document must be a ."  +  "docx" (spellcheck-flagged)  +  " Microsoft Word file"
#    This text appears twice (Volunteer + Product update sections); both are
#    rewritten identically.
# ---------------------------------------------------------------------------

$oldSentence = "The document must be a .docx Microsoft Word file"

function Split-DocxSentence {
    param($paragraph)

    $range = $paragraph.Range

    # Pull the paragraph's own rsid* attributes off its existing XML so the
    # rewritten paragraph keeps them unchanged (only the runs inside change).
    $srcXml = $range.WordOpenXML
    $rsidR = ""
    $rsidRDefault = ""
    $rsidP = ""
    if ($srcXml -match '<w:p\s[^>]*\bw:rsidR="([^"]*)"') { $rsidR = $matches[1] }
    if ($srcXml -match '<w:p\s[^>]*\bw:rsidRDefault="([^"]*)"') { $rsidRDefault = $matches[1] }
    if ($srcXml -match '<w:p\s[^>]*\bw:rsidP="([^"]*)"') { $rsidP = $matches[1] }

    $pAttrs = ""
    if ($rsidR -ne "") { $pAttrs += ' w:rsidR="' + $rsidR + '"' }
    if ($rsidRDefault -ne "") { $pAttrs += ' w:rsidRDefault="' + $rsidRDefault + '"' }
    if ($rsidP -ne "") { $pAttrs += ' w:rsidP="' + $rsidP + '"' }

    $newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"' + $pAttrs + '>' + `
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr>' + `
        '<w:r><w:t>The document must be a .</w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/>' + `
        '<w:r><w:t>docx</w:t></w:r>' + `
        '<w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:t xml:space="preserve"> Microsoft Word file</w:t></w:r>' + `
        '</w:p>'

    [void]$range.InsertXML($newXml)
}

foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd("`r") -eq $oldSentence) {
        Split-DocxSentence $para
    }
}

# ---------------------------------------------------------------------------
# 3) Remove the "Maintaining The System" Heading1 paragraph together with the
#    blank paragraph immediately before and after it (right before sectPr).
# ---------------------------------------------------------------------------

$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd("`r") -eq "Maintaining The System") {
        $target = $para
    }
}

if ($target -ne $null) {
    $prev = $target.Previous()
    $next = $target.Next()
    $delStart = $target.Range.Start
    $delEnd = $target.Range.End
    if ($prev -ne $null -and $prev.Range.Text.TrimEnd("`r") -eq "") {
        $delStart = $prev.Range.Start
    }
    if ($next -ne $null -and $next.Range.Text.TrimEnd("`r") -eq "") {
        $delEnd = $next.Range.End
    }
    $killRange = $d.Range($delStart, $delEnd)
    [void]$killRange.Delete()
}
